# The deck currently carries two theme parts:
#   ppt/theme/theme1.xml  -> "Integral" / "Red Violet" colour scheme (the
#                             theme actually wired to the one slide master
#                             that drives every slide's look)
#   ppt/theme/theme2.xml  -> the stock "Office Theme" (only wired to the
#                             notes master, so it is not visibly used)
#
# The authored edit swaps the contents of those two parts: the slide
# master (theme1.xml) becomes the plain "Office Theme" colours, and the
# previously-applied "Integral"/"Red Violet" palette moves into
# theme2.xml.
#
# This COM host exposes theme colour editing only through
# SlideMaster.ColorScheme.Colors(<index>).RGB, which patches the colour
# scheme that backs ppt/theme/theme1.xml (the live design). Each index
# lines up 1:1 with the clrScheme child order used by OOXML:
#   1 dk1   2 lt1   3 dk2   4 lt2
#   5 accent1  6 accent2  7 accent3  8 accent4  9 accent5  10 accent6
#   11 hlink   12 folHlink
#
# Re-apply the stock "Office Theme" palette (the colours presently sitting
# in theme2.xml) onto the live design so theme1.xml matches the target.

$p = $ppt.ActivePresentation
$design = $p.SlideMaster.ColorScheme

function Set-SchemeColor($index, $r, $g, $b) {
    $design.Colors($index).RGB = $r + ($g * 256) + ($b * 65536)
}

Set-SchemeColor 1  0x00 0x00 0x00   # dk1
Set-SchemeColor 2  0xFF 0xFF 0xFF   # lt1
Set-SchemeColor 3  0x44 0x54 0x6A   # dk2
Set-SchemeColor 4  0xE7 0xE6 0xE6   # lt2
Set-SchemeColor 5  0x5B 0x9B 0xD5   # accent1
Set-SchemeColor 6  0xED 0x7D 0x31   # accent2
Set-SchemeColor 7  0xA5 0xA5 0xA5   # accent3
Set-SchemeColor 8  0xFF 0xC0 0x00   # accent4
Set-SchemeColor 9  0x44 0x72 0xC4   # accent5
Set-SchemeColor 10 0x70 0xAD 0x47   # accent6
Set-SchemeColor 11 0x05 0x63 0xC1   # hlink
Set-SchemeColor 12 0x95 0x4F 0x72   # folHlink
